# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# zh-cn and de-de report sheets, as produced by a fresh handback report
# generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 09:07:20"
$wsZhCn.Range("H2").Value = "2016-03-22 09:07:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 09:07:24"
$wsDeDe.Range("H2").Value = "2016-03-22 09:07:49"
